# Apply the StructureDefinition-cobra-indicator.xlsx edits:
#  - Metadata sheet: URL / Version / Date / Publisher updates (Alvearie -> LinuxForHealth)
#  - Elements sheet: clear the ele-1/ext-1 constraint(s) text on the root Extension row

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/cobra-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Clear the Constraint(s) value for the root Extension element (row 2)
$elements.Range("AI2").Value = ""

# The Extension.url row repeats the canonical URL as its Fixed Value; keep it
# in sync with the updated URL on the Metadata sheet.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/cobra-indicator"
